# "scraping better with date" - append newly-scraped arrival rows (Monday,
# Jan 09) to the WMI_Arrivals worksheet, continuing the NUMBER sequence and
# DIFFERENCE calculation already used by the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: NUMBER, DATE, TIME, FLIGHT, FROM, SHORT, AIRLINE, MODEL,
#             AIRCFAT ID, STATUS, DIFFERENCE
$newRows = @(
    @(58, 'Monday, Jan 09', '2:30 PM', 'FR8011', 'Dublin', '(DUB)', 'Ryanair ', 'B38M', '(EI-HES)', '2:46 PM', '0 hours, 16 minutes'),
    @(59, 'Monday, Jan 09', '2:45 PM', 'FR4524', 'Edinburgh', '(EDI)', 'Ryanair ', 'B738', '(SP-RKP)', '3:03 PM', '0 hours, 18 minutes'),
    @(60, 'Monday, Jan 09', '3:15 PM', 'FR6944', 'Barcelona', '(BCN)', 'Ryanair ', 'B738', '(EI-DYC)', '2:59 PM', '0 hours, -16 minutes'),
    @(61, 'Monday, Jan 09', '3:20 PM', 'FR4534', 'Porto', '(OPO)', 'Ryanair ', 'B38M', '(SP-RZO)', '2:56 PM', '0 hours, -24 minutes'),
    @(62, 'Monday, Jan 09', '4:30 PM', 'FR4238', 'Bari', '(BRI)', 'Ryanair ', 'B738', '(SP-RSV)', '4:27 PM', '0 hours, -3 minutes'),
    @(63, 'Monday, Jan 09', '5:05 PM', 'FR1021', 'London', '(STN)', 'Ryanair ', 'B38M', '(SP-RZI)', '5:32 PM', '0 hours, 27 minutes'),
    @(64, 'Monday, Jan 09', '5:15 PM', 'FR1889', 'Paris', '(BVA)', 'Ryanair ', 'B738', '(SP-RKT)', '5:21 PM', '0 hours, 6 minutes'),
    @(65, 'Monday, Jan 09', '5:25 PM', 'FR1903', 'Eindhoven', '(EIN)', 'Ryanair ', 'B738', '(SP-RKL)', '5:16 PM', '0 hours, -9 minutes')
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data[0]   # A NUMBER
    $ws.Cells.Item($row, 2).Value = $data[1]   # B DATE
    $ws.Cells.Item($row, 3).Value = $data[2]   # C TIME
    $ws.Cells.Item($row, 4).Value = $data[3]   # D FLIGHT
    $ws.Cells.Item($row, 5).Value = $data[4]   # E FROM
    $ws.Cells.Item($row, 6).Value = $data[5]   # F SHORT
    $ws.Cells.Item($row, 7).Value = $data[6]   # G AIRLINE
    $ws.Cells.Item($row, 8).Value = $data[7]   # H MODEL
    $ws.Cells.Item($row, 9).Value = $data[8]   # I AIRCFAT ID
    $ws.Cells.Item($row, 10).Value = $data[9]  # J STATUS
    # K (DIFFERENCE header column) intentionally left blank, matching the
    # existing rows where the computed difference actually lives in L.
    $ws.Cells.Item($row, 12).Value = $data[10] # L DIFFERENCE
}
